$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.184.04'
$ws.Range("E2").Value = '  +2.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.434.57'
$ws.Range("E3").Value = '  +2.16%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '409.53'
$ws.Range("E5").Value = '  +0.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.76'
$ws.Range("E6").Value = '  -4.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.634'
$ws.Range("E7").Value = '  +7.41%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.755'
$ws.Range("E9").Value = '  +11.90%  '

$ws.Range("E10").Value = '  +18.66%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '43.09'
$ws.Range("E11").Value = '  +1.29%  '

$ws.Range("E12").Value = '  -0.31%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.88'
$ws.Range("E13").Value = '  +6.74%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.46'
$ws.Range("E14").Value = '  +4.53%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000194'
$ws.Range("E15").Value = '  +52.69%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.448.10'
$ws.Range("E16").Value = '  +2.68%  '

$ws.Range("E17").Value = '  +3.39%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.084.89'
$ws.Range("E18").Value = '  +1.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.36'
$ws.Range("E19").Value = '  +3.54%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '388.33'
$ws.Range("E20").Value = '  +25.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '88.62'
$ws.Range("E21").Value = '  +5.79%  '

$ws.Range("E22").Value = '  -0.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.36'
$ws.Range("E23").Value = '  +5.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.23'
$ws.Range("E24").Value = '  +3.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '32.05'
$ws.Range("E25").Value = '  +9.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.81'
$ws.Range("E26").Value = '  +0.59%  '

$ws.Range("E27").Value = '  +1.81%  '

$ws.Range("E28").Value = '  +3.49%  '

$ws.Range("E29").Value = '  +9.38%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '44.10'
$ws.Range("E30").Value = '  +7.16%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.116'
$ws.Range("E31").Value = '  -0.13%  '

$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.171'
$ws.Range("E32").Value = '  -0.65%  '

$ws.Range("E33").Value = '  +4.60%  '

$ws.Range("E34").Value = '  +0.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0491'
$ws.Range("E35").Value = '  +2.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '52.39'
$ws.Range("E36").Value = '  +1.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.37'
$ws.Range("E38").Value = '  -1.17%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.92'
$ws.Range("E39").Value = '  +0.72%  '

$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.133'
$ws.Range("E40").Value = '  +7.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.313'
$ws.Range("E41").Value = '  +9.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '140.88'
$ws.Range("E42").Value = '  +2.88%  '

$ws.Range("E43").Value = '  -0.44%  '

$ws.Range("E44").Value = '  -0.39%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.73'
$ws.Range("E45").Value = '  +0.83%  '

$ws.Range("E46").Value = '  +4.35%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.28'
$ws.Range("E47").Value = '  +4.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.119.87'
$ws.Range("E48").Value = '  +0.06%  '

$ws.Range("E49").Value = '  -0.20%  '

$ws.Range("E50").Value = '  +2.72%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0367'
$ws.Range("E51").Value = '  +6.72%  '
